$wb = $excel.ActiveWorkbook

# Add a new worksheet at the end of the workbook and name it for the 2025-12-24 week
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "magapoke_2025-12-24"

# Header row
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Ranking data (rank, title) for the 2025-12-24 week
$titles = @(
  "きゃわるり方程式",
  "黒月のイェルクナハト",
  "スルガメテオ",
  "ドリーム☆ジャンボ☆ガール",
  "アイドラトリィ",
  "お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！",
  "黄昏町プリズナーズ",
  "K-9~警視庁公安部公安第9課異能対策係~",
  "せいぶつ部の田辺くん",
  "篝家の８兄弟",
  "ゼロとヒャク",
  "ハードワーカー中田",
  "生きたがりの人狼",
  "ナキナギ",
  "ともだちづくり",
  "ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜",
  "普通の本はありません！",
  "明智ナンバーワン",
  "夜鐘のキト",
  "屋根の下のアルテミス",
  "ハナバス　苔石花江のバスケ論",
  "その青春",
  "MYS",
  "限界集落を脱村した錬金術士、都会で`"最強`"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～",
  "歪みの虜",
  "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～",
  "春くらり",
  "GURU",
  "ナマイキ旭ちゃんをわからせたい",
  "平成転生",
  "鳴るさんだぁ",
  "卒業アルバムの彼女たち",
  "じゅーくぼっくす",
  "JK Biker",
  "皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～",
  "君が監督！",
  "〈小市民〉 春期限定いちごタルト事件",
  "白銀のキュイジーヌ～明治外交官の料理人～",
  "追放されなかった男　～二度目の人生は土下座から始まりました～",
  "花子狩り",
  "ハプスブルク家の華麗なる受難",
  "永久のユウグレ",
  "人生逆転ダンジョン",
  "眠れる森のレガ",
  "イエティ、とある日々"
)

for ($i = 0; $i -lt $titles.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $i + 1
  $ws.Cells.Item($row, 2).Value = $titles[$i]
}

# Restore the originally active sheet (the new sheet becoming active on
# creation is just a side effect of adding it)
$wb.Worksheets.Item(1).Activate()
